# Actualizacion automatica del mapa (2025-12-04 10:12:21)
#
# This script reproduces, via Excel COM-interop calls, the edits described by
# the xml diff: two OT values filled in (E90/E91), two existing "NEW" rows
# (old rows 95/96) replaced in-place by two newly reported poste incidents
# (while the previously-there data is pushed down two rows), and three brand
# new rows appended at the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing TEXT storage, so that
# numeric-looking strings (leading zeros, leading "-", bare digits used as
# "Comuna" codes, trailing spaces, date-looking strings, etc.) are kept as
# literal text instead of being auto-coerced into numbers/dates by Excel.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

# ---------------------------------------------------------------------
# 1) Two OT numbers that were "Pendiente ADM" now have a real OT code.
# ---------------------------------------------------------------------
Set-TextValue $ws.Range("E90") "01675746 "
Set-TextValue $ws.Range("E91") "01639693 "

# ---------------------------------------------------------------------
# 2) Insert two fresh rows right before the old row 95, pushing the old
#    rows 95-97 down to 97-99 (their contents travel with them
#    automatically), then populate the two newly-opened rows (95 & 96)
#    with the newly reported incidents.
# ---------------------------------------------------------------------
$ws.Rows.Item(95).Insert()
$ws.Rows.Item(95).Insert()

# --- New row 95 ---
Set-TextValue $ws.Range("A95") "-688"
Set-TextValue $ws.Range("B95") "11/27/2025"
$ws.Range("C95").Value = "Murguiondo 4001"
Set-TextValue $ws.Range("D95") "8"
Set-TextValue $ws.Range("E95") "810862571"
$ws.Range("F95").Value = "NEW"
$ws.Range("G95").Value = "Pendiente"
$ws.Range("H95").Value = "corroida"
$ws.Range("I95").Value = 1
$ws.Range("J95").Value = "Cambio"
$ws.Range("K95").Value = "Sin equipos"
$ws.Range("L95").Value = "Pasante"
$ws.Range("M95").Value = -58.477556
$ws.Range("N95").Value = -34.675565
$ws.Range("O95").Value = "Boedo"
$ws.Range("P95").Value = "Capital Sur"
$ws.Range("Q95").Value = "PAV-V"
$ws.Range("R95").Value = "Fuera de Poligono OVL"

# --- New row 96 ---
Set-TextValue $ws.Range("A96") "-689"
Set-TextValue $ws.Range("B96") "11/27/2025"
$ws.Range("C96").Value = "Federico Garcia Lorca 285"
Set-TextValue $ws.Range("D96") "6"
Set-TextValue $ws.Range("E96") "810863096"
$ws.Range("F96").Value = "NEW"
$ws.Range("G96").Value = "Pendiente"
$ws.Range("H96").Value = "corroida"
$ws.Range("I96").Value = 1
$ws.Range("J96").Value = "Cambio"
$ws.Range("K96").Value = "Sin equipos"
$ws.Range("L96").Value = "Pasante"
$ws.Range("M96").Value = -58.445274
$ws.Range("N96").Value = -34.618595
$ws.Range("O96").Value = "Almagro"
$ws.Range("P96").Value = "Capital Sur"
$ws.Range("Q96").Value = "NRA-A"
$ws.Range("R96").Value = "Fuera de Poligono OVL"

# ---------------------------------------------------------------------
# 3) Append three brand-new incident rows (100, 101, 102) at the bottom.
# ---------------------------------------------------------------------

# --- Row 100 ---
Set-TextValue $ws.Range("A100") "-691"
Set-TextValue $ws.Range("B100") "12/1/2025"
$ws.Range("C100").Value = "Comodoro Rivadavia 1989"
Set-TextValue $ws.Range("D100") "13"
$ws.Range("E100").Value = "Pendiente ADM"
$ws.Range("F100").Value = "NEW"
$ws.Range("G100").Value = "Pendiente"
$ws.Range("H100").Value = "inclinada"
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = "Cambio"
$ws.Range("K100").Value = "Sin equipos"
$ws.Range("L100").Value = "Pasante"
$ws.Range("O100").Value = "No ubicado"
$ws.Range("P100").Value = "No clasificado, consultar con mantenimiento"
$ws.Range("Q100").Value = "No ubicado"
$ws.Range("R100").Value = "Fuera de Poligono OVL"

# --- Row 101 ---
Set-TextValue $ws.Range("A101") "-693"
Set-TextValue $ws.Range("B101") "12/1/2025"
$ws.Range("C101").Value = "Santos Tome 5015"
Set-TextValue $ws.Range("D101") "10"
$ws.Range("E101").Value = "Pendiente ADM"
$ws.Range("F101").Value = "NEW"
$ws.Range("G101").Value = "Pendiente"
$ws.Range("H101").Value = "base corroida correr 50cm por garaje"
$ws.Range("I101").Value = 1
$ws.Range("J101").Value = "Cambio"
$ws.Range("K101").Value = "Sin equipos"
$ws.Range("L101").Value = "Pasante"
$ws.Range("O101").Value = "No ubicado"
$ws.Range("P101").Value = "No clasificado, consultar con mantenimiento"
$ws.Range("Q101").Value = "No ubicado"
$ws.Range("R101").Value = "Fuera de Poligono OVL"

# --- Row 102 ---
Set-TextValue $ws.Range("A102") "-695"
Set-TextValue $ws.Range("B102") "12/3/2025"
$ws.Range("C102").Value = "Estomba 2626"
Set-TextValue $ws.Range("D102") "12"
$ws.Range("E102").Value = "Pendiente ADM"
$ws.Range("F102").Value = "NEW"
$ws.Range("G102").Value = "Pendiente"
$ws.Range("H102").Value = "desmontar"
$ws.Range("I102").Value = 1
$ws.Range("J102").Value = "Desmonte"
$ws.Range("K102").Value = "Sin equipos"
$ws.Range("L102").Value = "Pasante"
$ws.Range("M102").Value = -58.47538
$ws.Range("N102").Value = -34.566015
$ws.Range("O102").Value = "Colegiales"
$ws.Range("P102").Value = "Capital Norte"
$ws.Range("Q102").Value = "PUE-E"
$ws.Range("R102").Value = "Fuera de Poligono OVL"
